$d = $word.ActiveDocument

function Break-Text {
    param(
        [string]$FindText,
        [string]$ReplaceText
    )
    $rng = $d.Content
    $ok = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $FindText"
    }
}

# --- Paragraph 1: "Objetivos" body text -----------------------------------
Break-Text "cerâmicos.Preparar" "cerâmicos.^lPreparar"
Break-Text "cerâmicos.     Desenvolver" "cerâmicos.     ^lDesenvolver"

# --- Paragraph 2: "Programa" body text -------------------------------------
Break-Text "Introdução1.1A indústria cerâmica" "Introdução^l1.1A indústria cerâmica"
Break-Text "econômica e social 1.2Processamento" "econômica e social ^l1.2Processamento"
Break-Text "processos especiais1.3Produtos cerâmicos" "processos especiais^l1.3Produtos cerâmicos"
Break-Text "da aplicação e novos usos2Características" "da aplicação e novos usos^l^l2Características"
Break-Text "dos sólidos cerâmicos2.1Estruturas dos Cristais" "dos sólidos cerâmicos^l2.1Estruturas dos Cristais"
Break-Text "displasivas e reconstrutivas2.2Estruturas dos Vidros" "displasivas e reconstrutivas^l2.2Estruturas dos Vidros"
Break-Text "em sistemas vítreos, 2.3Imperfeições estruturais" "em sistemas vítreos, ^l2.3Imperfeições estruturais"
Break-Text "não estequiométricos e discordâncias.2.4Superfície" "não estequiométricos e discordâncias.^l2.4Superfície"
Break-Text "molhabilidade e distribuição de fases. 2.5Mobilidade Atômica" "molhabilidade e distribuição de fases. ^l2.5Mobilidade Atômica"
Break-Text "contorno de grão e superfície e difusão em vidros.3Propriedades" "contorno de grão e superfície e difusão em vidros.^l^l3Propriedades"
Break-Text "3Propriedades de Cerâmicas:3.1Propriedades Térmicas" "3Propriedades de Cerâmicas:^l3.1Propriedades Térmicas"
Break-Text "brilho de superfície, aplicações3.2Propriedades Ópticas" "brilho de superfície, aplicações^l3.2Propriedades Ópticas"
Break-Text "translucidez e opacidade, absorção e cor, 3.3Deformação Plástica" "translucidez e opacidade, absorção e cor, ^l3.3Deformação Plástica"
Break-Text "fluxo viscoso em líquidos e vidros.3.4Elasticidade" "fluxo viscoso em líquidos e vidros.^l3.4Elasticidade"
Break-Text "fratura por fluência e efeitos da microestrutura3.5Tensões térmicas" "fratura por fluência e efeitos da microestrutura^l3.5Tensões térmicas"

# --- Paragraph 3: "Bibliografia" body text ----------------------------------
Break-Text "John Wiley, c1976;2 - R.W. CAHN" "John Wiley, c1976;^l2 - R.W. CAHN"
Break-Text "Wiley-VCH, c2005;3 - BERGERON" "Wiley-VCH, c2005;^l3 - BERGERON"
Break-Text "Ceramic Society, 1984;4 - BROOK" "Ceramic Society, 1984;^l4 - BROOK"
Break-Text "Weinheim: VCH, 1996;5 - LEVIN" "Weinheim: VCH, 1996;^l5 - LEVIN"
Break-Text "American Ceramic, 1964;6 - R.C. BRADT" "American Ceramic, 1964;^l6 - R.C. BRADT"
Break-Text "and thermal shock.. New York: Plenum, 1996.7 - R.C. BRADT" "and thermal shock.. New York: Plenum, 1996.^l7 - R.C. BRADT"
Break-Text "and high-temperature behaviour.. New York: Plenum, 1996;8 - REED" "and high-temperature behaviour.. New York: Plenum, 1996;^l8 - REED"
Break-Text "New YorK: John Wiley, 1995;9 - TOMPSON" "New YorK: John Wiley, 1995;^l9 - TOMPSON"
Break-Text "Institute of Materials, 1993;10 - BARSOUM" "Institute of Materials, 1993;^l10 - BARSOUM"
Break-Text "The McGraw-Hill, 1997;11 - CHIANG" "The McGraw-Hill, 1997;^l11 - CHIANG"
Break-Text "John Wiley, 1997;12 - MENCIK" "John Wiley, 1997;^l12 - MENCIK"

Write-Output "Done"
